$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")
$ws.Activate()

# --- Record progress for the tasks worked on during the sprint ---
# Task 1 (rows 5-7), group assigned to Vanja Cvetković, marked Done
$ws.Range("D5").Value = "Vanja Cvetković"
$ws.Range("F5").Value = "Done"
$ws.Range("H5").Value = 3

$ws.Range("D6").Value = "Vanja Cvetković"
$ws.Range("F6").Value = "Done"
$ws.Range("J6").Value = 8

$ws.Range("F7").Value = "Done"
$ws.Range("J7").Value = 2

# Task 2 (rows 8-10), group assigned to Vanja Cvetković, marked Done
$ws.Range("D8").Value = "Vanja Cvetković"
$ws.Range("F8").Value = "Done"
$ws.Range("K8").Value = 3

$ws.Range("D9").Value = "Vanja Cvetković"
$ws.Range("F9").Value = "Done"
$ws.Range("M9").Value = 8

$ws.Range("F10").Value = "Done"
$ws.Range("M10").Value = 2

# Task 3 (rows 11-14)
$ws.Range("F11").Value = "Done"
$ws.Range("M11").Value = 1

$ws.Range("D14").Value = "Predrag Dimitrijević"
$ws.Range("F14").Value = "Done"
$ws.Range("J14").Value = 1

# Task 4 (rows 15-17)
$ws.Range("L15").Value = 2

# --- Adjust conditional formatting ranges to include row 11 in the first group ---
$fcs2 = $ws.Range("F11:F14").FormatConditions
$fcs2.Item(1).ModifyAppliesToRange($ws.Range("F12:F14"))

$fcs1 = $ws.Range("F5:F10").FormatConditions
$fcs1.Item(1).ModifyAppliesToRange($ws.Range("F5:F11"))

# --- Update the selection on the Sprint sheet ---
$ws.Range("L16").Select()
